# Scheduled data refresh: update cached market-board figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# across the per-crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1321.1538
$ws.Cells.Item(18, 9).Value = 1145.6522
$ws.Cells.Item(18, 10).Value = 2666.6667
$ws.Cells.Item(18, 11).Value = 1145.6522
$ws.Cells.Item(18, 12).Value = 2666.6667
$ws.Cells.Item(18, 13).Value = -861.6522
$ws.Cells.Item(18, 14).Value = -3234.6667
$ws.Cells.Item(19, 8).Value = 2486.1428
$ws.Cells.Item(19, 9).Value = 5269
$ws.Cells.Item(19, 10).Value = 1167.9474
$ws.Cells.Item(19, 11).Value = 5269
$ws.Cells.Item(19, 12).Value = 1167.9474
$ws.Cells.Item(19, 13).Value = -5094
$ws.Cells.Item(19, 14).Value = -1517.9474
$ws.Cells.Item(62, 8).Value = 3010.8125
$ws.Cells.Item(62, 9).Value = 3061.182
$ws.Cells.Item(62, 10).Value = 2900
$ws.Cells.Item(62, 11).Value = 3061.182
$ws.Cells.Item(62, 12).Value = 2900
$ws.Cells.Item(62, 13).Value = -2437.182
$ws.Cells.Item(62, 14).Value = -4148
$ws.Cells.Item(65, 8).Value = 3010.8125
$ws.Cells.Item(65, 9).Value = 3061.182
$ws.Cells.Item(65, 10).Value = 2900
$ws.Cells.Item(65, 11).Value = 15305.91
$ws.Cells.Item(65, 12).Value = 14500
$ws.Cells.Item(65, 13).Value = -12185.91
$ws.Cells.Item(65, 14).Value = -20740
$ws.Cells.Item(116, 8).Value = 1783.3334
$ws.Cells.Item(116, 9).Value = 1500
$ws.Cells.Item(116, 11).Value = 1500
$ws.Cells.Item(116, 13).Value = 1942
$ws.Cells.Item(132, 8).Value = 3175710.5
$ws.Cells.Item(132, 9).Value = 3247817.5
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 9743452.5
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -9740922.5
$ws.Cells.Item(132, 14).Value = -14060
$ws.Cells.Item(137, 8).Value = 1644.7667
$ws.Cells.Item(137, 9).Value = 1387.15
$ws.Cells.Item(137, 10).Value = 2160
$ws.Cells.Item(137, 11).Value = 4161.450000000001
$ws.Cells.Item(137, 12).Value = 6480
$ws.Cells.Item(137, 13).Value = -1611.450000000001
$ws.Cells.Item(137, 14).Value = -11580

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1025.9333
$ws.Cells.Item(45, 9).Value = 952.1579
$ws.Cells.Item(45, 10).Value = 1153.3636
$ws.Cells.Item(45, 11).Value = 952.1579
$ws.Cells.Item(45, 12).Value = 1153.3636
$ws.Cells.Item(45, 13).Value = -575.1579
$ws.Cells.Item(45, 14).Value = -1907.3636
$ws.Cells.Item(61, 8).Value = 1825.5238
$ws.Cells.Item(61, 9).Value = 955.73334
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 955.73334
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = -743.73334
$ws.Cells.Item(61, 14).Value = -4424
$ws.Cells.Item(63, 8).Value = 836083.3
$ws.Cells.Item(63, 9).Value = 1002100
$ws.Cells.Item(63, 10).Value = 6000
$ws.Cells.Item(63, 11).Value = 1002100
$ws.Cells.Item(63, 12).Value = 6000
$ws.Cells.Item(63, 13).Value = -1001414
$ws.Cells.Item(63, 14).Value = -7372
$ws.Cells.Item(66, 8).Value = 836083.3
$ws.Cells.Item(66, 9).Value = 1002100
$ws.Cells.Item(66, 10).Value = 6000
$ws.Cells.Item(66, 11).Value = 5010500
$ws.Cells.Item(66, 12).Value = 30000
$ws.Cells.Item(66, 13).Value = -5007068
$ws.Cells.Item(66, 14).Value = -36864
$ws.Cells.Item(74, 8).Value = 1735.6
$ws.Cells.Item(74, 9).Value = 1453.6
$ws.Cells.Item(74, 11).Value = 1453.6
$ws.Cells.Item(74, 13).Value = -579.5999999999999
$ws.Cells.Item(77, 8).Value = 1735.6
$ws.Cells.Item(77, 9).Value = 1453.6
$ws.Cells.Item(77, 11).Value = 7268
$ws.Cells.Item(77, 13).Value = -2900
$ws.Cells.Item(102, 8).Value = 2160
$ws.Cells.Item(102, 9).Value = 2160
$ws.Cells.Item(102, 11).Value = 2160
$ws.Cells.Item(102, 13).Value = -538
$ws.Cells.Item(110, 8).Value = 1356.55
$ws.Cells.Item(110, 9).Value = 951.2857
$ws.Cells.Item(110, 10).Value = 2302.1667
$ws.Cells.Item(110, 11).Value = 951.2857
$ws.Cells.Item(110, 12).Value = 2302.1667
$ws.Cells.Item(110, 13).Value = 1093.7143
$ws.Cells.Item(110, 14).Value = -6392.1667
$ws.Cells.Item(122, 8).Value = 1986.2667
$ws.Cells.Item(122, 9).Value = 1240.2
$ws.Cells.Item(122, 10).Value = 2359.3
$ws.Cells.Item(122, 11).Value = 3720.6
$ws.Cells.Item(122, 12).Value = 7077.900000000001
$ws.Cells.Item(122, 13).Value = -1270.6
$ws.Cells.Item(122, 14).Value = -11977.9
$ws.Cells.Item(136, 8).Value = 1825.5238
$ws.Cells.Item(136, 9).Value = 955.73334
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 2867.20002
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -317.2000200000002
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1683.4348
$ws.Cells.Item(105, 9).Value = 1362.375
$ws.Cells.Item(105, 10).Value = 2417.2856
$ws.Cells.Item(105, 11).Value = 1362.375
$ws.Cells.Item(105, 12).Value = 2417.2856
$ws.Cells.Item(105, 13).Value = 384.625
$ws.Cells.Item(105, 14).Value = -5911.2856
$ws.Cells.Item(107, 8).Value = 794.1053000000001
$ws.Cells.Item(107, 9).Value = 794.5625
$ws.Cells.Item(107, 10).Value = 791.6667
$ws.Cells.Item(107, 11).Value = 794.5625
$ws.Cells.Item(107, 12).Value = 791.6667
$ws.Cells.Item(107, 13).Value = 1125.4375
$ws.Cells.Item(107, 14).Value = -4631.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 965.55554
$ws.Cells.Item(16, 10).Value = 931.6667
$ws.Cells.Item(16, 12).Value = 931.6667
$ws.Cells.Item(16, 14).Value = -1505.6667
$ws.Cells.Item(22, 8).Value = 274.17648
$ws.Cells.Item(22, 9).Value = 300.7857
$ws.Cells.Item(22, 11).Value = 300.7857
$ws.Cells.Item(22, 13).Value = 49.21429999999998
$ws.Cells.Item(31, 8).Value = 5716553.5
$ws.Cells.Item(31, 9).Value = 2334.2354
$ws.Cells.Item(31, 10).Value = 200000000
$ws.Cells.Item(31, 11).Value = 2334.2354
$ws.Cells.Item(31, 12).Value = 200000000
$ws.Cells.Item(31, 13).Value = -2039.2354
$ws.Cells.Item(31, 14).Value = -200000590
$ws.Cells.Item(34, 8).Value = 5716553.5
$ws.Cells.Item(34, 9).Value = 2334.2354
$ws.Cells.Item(34, 10).Value = 200000000
$ws.Cells.Item(34, 11).Value = 2334.2354
$ws.Cells.Item(34, 12).Value = 200000000
$ws.Cells.Item(34, 13).Value = -2132.2354
$ws.Cells.Item(34, 14).Value = -200000404
$ws.Cells.Item(86, 8).Value = 1998.3
$ws.Cells.Item(86, 9).Value = 1861.4
$ws.Cells.Item(86, 10).Value = 2135.2
$ws.Cells.Item(86, 11).Value = 1861.4
$ws.Cells.Item(86, 12).Value = 2135.2
$ws.Cells.Item(86, 13).Value = -738.4000000000001
$ws.Cells.Item(86, 14).Value = -4381.2
$ws.Cells.Item(89, 8).Value = 1998.3
$ws.Cells.Item(89, 9).Value = 1861.4
$ws.Cells.Item(89, 10).Value = 2135.2
$ws.Cells.Item(89, 11).Value = 9307
$ws.Cells.Item(89, 12).Value = 10676
$ws.Cells.Item(89, 13).Value = -3691
$ws.Cells.Item(89, 14).Value = -21908
$ws.Cells.Item(105, 8).Value = 733.94446
$ws.Cells.Item(105, 9).Value = 703.75
$ws.Cells.Item(105, 10).Value = 758.1
$ws.Cells.Item(105, 11).Value = 703.75
$ws.Cells.Item(105, 12).Value = 758.1
$ws.Cells.Item(105, 13).Value = 1043.25
$ws.Cells.Item(105, 14).Value = -4252.1
$ws.Cells.Item(107, 8).Value = 548.4400000000001
$ws.Cells.Item(107, 10).Value = 622.75
$ws.Cells.Item(107, 12).Value = 622.75
$ws.Cells.Item(107, 14).Value = -4462.75
$ws.Cells.Item(113, 8).Value = 965.55554
$ws.Cells.Item(113, 10).Value = 931.6667
$ws.Cells.Item(113, 12).Value = 931.6667
$ws.Cells.Item(113, 14).Value = -5271.6667
$ws.Cells.Item(132, 8).Value = 1871.35
$ws.Cells.Item(132, 9).Value = 1195.6666
$ws.Cells.Item(132, 11).Value = 3586.9998
$ws.Cells.Item(132, 13).Value = -1056.9998
$ws.Cells.Item(134, 8).Value = 834.3036
$ws.Cells.Item(134, 9).Value = 760.2708
$ws.Cells.Item(134, 10).Value = 1278.5
$ws.Cells.Item(134, 11).Value = 2280.8124
$ws.Cells.Item(134, 12).Value = 3835.5
$ws.Cells.Item(134, 13).Value = 254.1876000000002
$ws.Cells.Item(134, 14).Value = -8905.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1468.3158
$ws.Cells.Item(122, 9).Value = 1509.1818
$ws.Cells.Item(122, 10).Value = 1412.125
$ws.Cells.Item(122, 11).Value = 13582.6362
$ws.Cells.Item(122, 12).Value = 12709.125
$ws.Cells.Item(122, 13).Value = -11132.6362
$ws.Cells.Item(122, 14).Value = -17609.125
$ws.Cells.Item(131, 8).Value = 3301276.5
$ws.Cells.Item(131, 9).Value = 9569.916999999999
$ws.Cells.Item(131, 10).Value = 5096752.5
$ws.Cells.Item(131, 11).Value = 28709.751
$ws.Cells.Item(131, 12).Value = 15290257.5
$ws.Cells.Item(131, 13).Value = -23669.751
$ws.Cells.Item(131, 14).Value = -15300337.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 142863090
$ws.Cells.Item(92, 10).Value = 142863090
$ws.Cells.Item(92, 12).Value = 142863090
$ws.Cells.Item(92, 14).Value = -142866834
$ws.Cells.Item(113, 8).Value = 50001600
$ws.Cells.Item(113, 9).Value = 83334340
$ws.Cells.Item(113, 11).Value = 83334340
$ws.Cells.Item(113, 13).Value = -83332170
$ws.Cells.Item(132, 8).Value = 67006.55
$ws.Cells.Item(132, 10).Value = 3999.6667
$ws.Cells.Item(132, 12).Value = 11999.0001
$ws.Cells.Item(132, 14).Value = -17059.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1702.5
$ws.Cells.Item(40, 9).Value = 1700
$ws.Cells.Item(40, 10).Value = 1708.75
$ws.Cells.Item(40, 11).Value = 1700
$ws.Cells.Item(40, 12).Value = 1708.75
$ws.Cells.Item(40, 13).Value = -1564
$ws.Cells.Item(40, 14).Value = -1980.75
$ws.Cells.Item(82, 8).Value = 1067.6471
$ws.Cells.Item(82, 9).Value = 1004.5455
$ws.Cells.Item(82, 10).Value = 1183.3334
$ws.Cells.Item(82, 11).Value = 1004.5455
$ws.Cells.Item(82, 12).Value = 1183.3334
$ws.Cells.Item(82, 13).Value = -643.5454999999999
$ws.Cells.Item(82, 14).Value = -1905.3334
$ws.Cells.Item(85, 8).Value = 1067.6471
$ws.Cells.Item(85, 9).Value = 1004.5455
$ws.Cells.Item(85, 10).Value = 1183.3334
$ws.Cells.Item(85, 11).Value = 1004.5455
$ws.Cells.Item(85, 12).Value = 1183.3334
$ws.Cells.Item(85, 13).Value = 243.4545000000001
$ws.Cells.Item(85, 14).Value = -3679.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 50275
$ws.Cells.Item(92, 10).Value = 50275
$ws.Cells.Item(92, 12).Value = 50275
$ws.Cells.Item(92, 14).Value = -55267
$ws.Cells.Item(94, 8).Value = 55000
$ws.Cells.Item(94, 10).Value = 55000
$ws.Cells.Item(94, 12).Value = 55000
$ws.Cells.Item(94, 14).Value = -56802
$ws.Cells.Item(122, 8).Value = 1137.2273
$ws.Cells.Item(122, 9).Value = 982.4375
$ws.Cells.Item(122, 10).Value = 1550
$ws.Cells.Item(122, 11).Value = 2947.3125
$ws.Cells.Item(122, 12).Value = 4650
$ws.Cells.Item(122, 13).Value = -497.3125
$ws.Cells.Item(122, 14).Value = -9550
$ws.Cells.Item(126, 8).Value = 5801.5
$ws.Cells.Item(126, 9).Value = 7048.75
$ws.Cells.Item(126, 10).Value = 812.5
$ws.Cells.Item(126, 11).Value = 21146.25
$ws.Cells.Item(126, 12).Value = 2437.5
$ws.Cells.Item(126, 13).Value = -18676.25
$ws.Cells.Item(126, 14).Value = -7377.5
